# Generate Report for Handback
# Mirrors the localization handback report generator: for each localized
# sheet (zh-cn, de-de), fill in the "Latest Target File" / "Latest Handback
# File" / "Latest Handback DateTime" columns for the two rows that just
# completed handback, flip their Status from "In Translation" to
# "Handed back: in sync with en-US", and hyperlink the newly filled target
# file cell the same way the source-file cell already is.

$wb = $excel.ActiveWorkbook

$githubBase = "https://github.com/OpenLocalizationTestOrg/oltest/blob/938b017f3e30ef2538de3452342a17cddf9694a8/e2e/"
$mdFile1 = "2c10ed6d-07f2-473a-a344-52cf11ef3295.md"
$mdFile2 = "38533bbb-9720-4552-bfce-4075550e5afe.md"

$statusText = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Status column -> handed back
$wsZh.Range("C2:C3").Value = $statusText

# Latest Target File (I) / Latest Handback File (J) / Latest Handback DateTime (K)
$wsZh.Range("J2").Value = "2c10ed6d-07f2-473a-a344-52cf11ef3295.e31e422a680ec9c7ed7e9d8daaabf4b899344e4b.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-13 06:29:59"

$wsZh.Range("J3").Value = "38533bbb-9720-4552-bfce-4075550e5afe.9c36f8a28459f9632c53380567221c8e65ce1823.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-08-13 06:29:59"

# Rebuild the hyperlinks collection so the new "Latest Target File" links sit
# next to the existing "File Name" links, same as the source workbook does.
$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $githubBase + $mdFile1, "", "", $mdFile1)
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $githubBase + $mdFile1, "", "", $mdFile1)
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $githubBase + $mdFile2, "", "", $mdFile2)
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $githubBase + $mdFile2, "", "", $mdFile2)

# Column widths widen to fit the longer status text / new file-name columns.
$wsZh.Range("C1").ColumnWidth = 29.15
$wsZh.Range("I1").ColumnWidth = 39.1
$wsZh.Range("J1").ColumnWidth = 39.1

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# Status column -> handed back
$wsDe.Range("C2:C3").Value = $statusText

# Latest Target File (I) / Latest Handback File (J) / Latest Handback DateTime (K)
$wsDe.Range("J2").Value = "2c10ed6d-07f2-473a-a344-52cf11ef3295.e31e422a680ec9c7ed7e9d8daaabf4b899344e4b.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-13 06:30:16"

$wsDe.Range("J3").Value = "38533bbb-9720-4552-bfce-4075550e5afe.9c36f8a28459f9632c53380567221c8e65ce1823.de-de.xlf"
$wsDe.Range("K3").Value = "2016-08-13 06:30:16"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $githubBase + $mdFile1, "", "", $mdFile1)
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $githubBase + $mdFile1, "", "", $mdFile1)
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $githubBase + $mdFile2, "", "", $mdFile2)
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $githubBase + $mdFile2, "", "", $mdFile2)

$wsDe.Range("C1").ColumnWidth = 29.15
$wsDe.Range("I1").ColumnWidth = 39.1
$wsDe.Range("J1").ColumnWidth = 39.1

# ---------------------------------------------------------------------
# Overview sheet - no cell values change, just widen the zh-cn/de-de status
# columns (E, F) to match the longer status text now shown there.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E1").ColumnWidth = 29.15
$wsOverview.Range("F1").ColumnWidth = 29.15
